$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.502.78"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "2.635.87"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'112.30"
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("D6").Value = "'326.03"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("E7").Value = "  -1.22%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.550"
$ws.Range("E9").Value = "  -0.83%  "
$ws.Range("D10").Value = "'39.60"
$ws.Range("E10").Value = "  -3.47%  "
$ws.Range("D11").Value = "'19.93"
$ws.Range("E11").Value = "  -1.00%  "
$ws.Range("D12").Value = "'0.0814"
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("E13").Value = "  +1.75%  "
$ws.Range("D14").Value = "'7.61"
$ws.Range("E14").Value = "  +3.39%  "
$ws.Range("D15").Value = "3.046.13"
$ws.Range("E15").Value = "  -0.50%  "
$ws.Range("D16").Value = "2.624.12"
$ws.Range("E16").Value = "  -1.14%  "
$ws.Range("D17").Value = "'0.855"
$ws.Range("E17").Value = "  -2.01%  "
$ws.Range("D18").Value = "49.446.01"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").Value = "'13.38"
$ws.Range("E19").Value = "  +2.01%  "
$ws.Range("E20").Value = "  -0.75%  "
$ws.Range("E21").Value = "  -1.89%  "
$ws.Range("E22").Value = "  -0.96%  "
$ws.Range("D23").Value = "'268.58"
$ws.Range("E23").Value = "  -3.14%  "
$ws.Range("D24").Value = "'69.14"
$ws.Range("E24").Value = "  -4.03%  "
$ws.Range("E25").Value = "  -0.90%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "'26.03"
$ws.Range("E27").Value = "  -2.79%  "
$ws.Range("D28").Value = "'10.14"
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("D30").Value = "'0.138"
$ws.Range("E30").Value = "  -3.17%  "
$ws.Range("D31").Value = "'34.46"
$ws.Range("E31").Value = "  -4.52%  "
$ws.Range("D32").Value = "'49.64"
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("E33").Value = "  +0.60%  "
$ws.Range("D34").Value = "'0.0814"
$ws.Range("E34").Value = "  +0.48%  "
$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").Value = "'19.12"
$ws.Range("E35").Value = "  -2.31%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("E37").Value = "  +2.32%  "
$ws.Range("D38").Value = "'2.02"
$ws.Range("E38").Value = "  -2.31%  "
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("D40").Value = "'129.56"
$ws.Range("E40").Value = "  +4.35%  "
$ws.Range("E41").Value = "  +6.00%  "
$ws.Range("D42").Value = "'23.25"
$ws.Range("E42").Value = "  +5.29%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0337"
$ws.Range("E43").Value = "  +6.88%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").Value = "'0.111"
$ws.Range("E44").Value = "  -0.97%  "
$ws.Range("D45").Value = "2.062.78"
$ws.Range("E45").Value = "  -0.79%  "
$ws.Range("E46").Value = "  -0.47%  "
$ws.Range("E47").Value = "  +7.59%  "
$ws.Range("E48").Value = "  -7.27%  "
$ws.Range("E49").Value = "  -2.78%  "
$ws.Range("D50").Value = "'5.22"
$ws.Range("E50").Value = "  -3.14%  "
$ws.Range("D51").Value = "'58.57"
$ws.Range("E51").Value = "  -1.57%  "
